# Applies "hybrid bold + color" highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) across the achievements /
# work-experience bullet points, matching the target OOXML diff.
#
# Word's Font.Color expects a BGR-packed long (0x00BBGGRR), so RGB 2C3E50
# is written as the hex literal 0x503E2C below.

$d = $word.ActiveDocument
$highlightColor = 0x503E2C

function Format-Metric($paragraphRange, [string]$metricText) {
    $found = $paragraphRange.Duplicate
    # MatchCase=$true, MatchWildcards=$false - literal substring search,
    # scoped to the paragraph's own Range so it can't bleed into other
    # paragraphs that happen to contain the same digits elsewhere.
    $ok = $found.Find.Execute($metricText, $true, $false, $false, $false, $false, $false, 1, $false, "", 0)
    if ($ok) {
        $found.Bold = 1
        $found.Font.Color = $highlightColor
    }
}

# Map of paragraph index -> list of metric substrings to bold+color,
# in left-to-right order as they occur in the paragraph text.
$targets = @(
    @{ Index = 10; Metrics = @("23%", "64%") },
    @{ Index = 12; Metrics = @("±4.2%", "±2.1%", "71%", "87%") },
    @{ Index = 13; Metrics = @("73.5%", "$4.7M") },
    @{ Index = 14; Metrics = @("$2") },
    @{ Index = 24; Metrics = @("57%") },
    @{ Index = 50; Metrics = @("$4.9M") },
    @{ Index = 51; Metrics = @("23%") },
    @{ Index = 53; Metrics = @("12,847") }
)

foreach ($entry in $targets) {
    $para = $d.Paragraphs.Item($entry.Index).Range
    foreach ($metric in $entry.Metrics) {
        Format-Metric $para $metric
    }
}
